$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 21:52"

# Row 22 - Asturias
$ws.Range("B22").Value = 1737
$ws.Range("C22").Value = 340
$ws.Range("D22").Value = 1287
$ws.Range("E22").Value = 110

# Row 31 - Murcia
$ws.Range("B31").Value = 1356
$ws.Range("C31").Value = 219
$ws.Range("D31").Value = 1049
$ws.Range("E31").Value = 88

# Row 54 - Melilla
$ws.Range("C54").Value = 16
$ws.Range("D54").Value = 75

# Row 55 - Ceuta
$ws.Range("C55").Value = 8
$ws.Range("D55").Value = 72
